$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.52586589497228919
$ws.Range("BO1").Value = 0.82934552474692558
$ws.Range("BP1").Value = 0.94461392975225622
$ws.Range("A2").Value = 0.82063710708296711
$ws.Range("C2").Value = 0.92775983118115701
$ws.Range("AM2").Value = 0.8325217656496352
$ws.Range("BP2").Value = 0.74164300771332559
$ws.Range("D3").Value = 0.71159669988627239
$ws.Range("E3").Value = 0.68493004877810648
$ws.Range("B4").Value = 0.97275154077556403
$ws.Range("F4").Value = 0.97683603788120854
$ws.Range("M4").Value = 0.76967993289726522
$ws.Range("F5").Value = 0.8089553199050209
$ws.Range("E7").Value = 0.6745868267956483
$ws.Range("F7").Value = 0.85770940888204339
$ws.Range("F8").Value = 0.81346494177800144
$ws.Range("G8").Value = 0.74938790112678522
$ws.Range("H9").Value = 0.86179024882318322
$ws.Range("AR9").Value = 0.99137398669302046
$ws.Range("G10").Value = 0.88883443733817147
$ws.Range("H10").Value = 0.75101290555419986
$ws.Range("J11").Value = 0.63860975164787881
$ws.Range("M11").Value = 0.58934580936657155
$ws.Range("K12").Value = 0.97947201057669853
$ws.Range("N12").Value = 0.98286168285439335
$ws.Range("AB12").Value = 0.69832973376718344
$ws.Range("N13").Value = 0.69092717423889227
$ws.Range("AR13").Value = 0.61081822172336797
$ws.Range("AA14").Value = 0.87419297277718355
$ws.Range("M15").Value = 0.67874066689628099
$ws.Range("N15").Value = 0.77660396699164602
$ws.Range("N16").Value = 0.97610898203238794
$ws.Range("O16").Value = 0.86115414069652485
$ws.Range("O17").Value = 0.70624884069963167
$ws.Range("R17").Value = 0.92142121137273647
$ws.Range("S17").Value = 0.70453088779627981
$ws.Range("P18").Value = 0.60079344618086472
$ws.Range("T18").Value = 0.7214539701728504
$ws.Range("BF18").Value = 0.94083519242447111
$ws.Range("K19").Value = 0.98600081123979755
$ws.Range("S20").Value = 0.87995389125853318
$ws.Range("V20").Value = 0.93694151777501111
$ws.Range("AP20").Value = 0.63026606150814179
$ws.Range("U23").Value = 0.58486508284247574
$ws.Range("Y23").Value = 0.81484581774960252
$ws.Range("V24").Value = 0.67548442552873822
$ws.Range("W24").Value = 0.99088897281603416
$ws.Range("Z24").Value = 0.87797880807668238
$ws.Range("BF24").Value = 0.82820164445782485
$ws.Range("V25").Value = 0.93132399293328205
$ws.Range("BL25").Value = 0.84567163131251055
$ws.Range("Y26").Value = 0.99442941775102123
$ws.Range("AB26").Value = 0.54468989813705582
$ws.Range("BA26").Value = 0.88304319282303778
$ws.Range("AC27").Value = 0.68204820536924271
$ws.Range("AD28").Value = 0.94677521819551269
$ws.Range("AD29").Value = 0.76414286115066021
$ws.Range("AE29").Value = 0.65830515056912398
$ws.Range("G30").Value = 0.82383135924426354
$ws.Range("P30").Value = 0.88846600285033694
$ws.Range("V30").Value = 0.77204087033858593
$ws.Range("AF31").Value = 0.88036423007468123
$ws.Range("AG32").Value = 0.81803068659764944
$ws.Range("J33").Value = 0.7065865316417872
$ws.Range("AE33").Value = 0.76511515250511697
$ws.Range("AH33").Value = 0.99301150901232416
$ws.Range("AY33").Value = 0.87828229944835234
$ws.Range("AF34").Value = 0.93050600629039826
$ws.Range("AJ34").Value = 0.92941108018328378
$ws.Range("U35").Value = 0.78646500446805256
$ws.Range("AG35").Value = 0.62233462126217987
$ws.Range("AH35").Value = 0.98360422023897365
$ws.Range("AJ35").Value = 0.79477939188450553
$ws.Range("AK35").Value = 0.75675782911688272
$ws.Range("AJ37").Value = 0.6824246299086274
$ws.Range("AL37").Value = 0.80911648352793031
$ws.Range("AM37").Value = 0.94036179976034084
$ws.Range("AJ38").Value = 0.57647156567530533
$ws.Range("AL39").Value = 0.91868043096812357
$ws.Range("AO39").Value = 0.77269094077638356
$ws.Range("D40").Value = 0.82174663774151369
$ws.Range("AL40").Value = 0.62622108465807891
$ws.Range("AP41").Value = 0.7903013302556563
$ws.Range("AQ41").Value = 0.93138871273436363
$ws.Range("AF42").Value = 0.74572321678430264
$ws.Range("AQ42").Value = 0.93215081316202353
$ws.Range("AR43").Value = 0.6493672676494826
$ws.Range("AS43").Value = 0.98671771648736495
$ws.Range("AS44").Value = 0.87310497443054058
$ws.Range("AU45").Value = 0.85152838435133682
$ws.Range("BG45").Value = 0.7565772300078859
$ws.Range("AU46").Value = 0.98569216765438927
$ws.Range("AT48").Value = 0.96355851260891878
$ws.Range("AU48").Value = 0.79171742818801993
$ws.Range("AU49").Value = 0.61259945124007886
$ws.Range("AV49").Value = 0.84664180437618741
$ws.Range("AY49").Value = 0.777946998776886
$ws.Range("BF49").Value = 0.56396815870029027
$ws.Range("AY50").Value = 0.99673472613799241
$ws.Range("L51").Value = 0.87989783535592059
$ws.Range("AZ51").Value = 0.77905651500637374
$ws.Range("AX52").Value = 0.75270017829148372
$ws.Range("BA52").Value = 0.82358678533522389
$ws.Range("AN53").Value = 0.92834326201941186
$ws.Range("BC53").Value = 0.75643235941666753
$ws.Range("AZ54").Value = 0.87111515672734785
$ws.Range("BA54").Value = 0.91683591051843716
$ws.Range("BD55").Value = 0.86219639242343149
$ws.Range("AI56").Value = 0.97690625993327485
$ws.Range("BB56").Value = 0.62360247021355042
$ws.Range("BC57").Value = 0.89218926164786572
$ws.Range("BB58").Value = 0.99506267482061239
$ws.Range("BD58").Value = 0.71424344293835751
$ws.Range("BE58").Value = 0.98188137681663601
$ws.Range("BG58").Value = 0.65391043002477067
$ws.Range("BH58").Value = 0.89091182924704282
$ws.Range("BE59").Value = 0.90947752109796798
$ws.Range("BI59").Value = 0.9448457170399116
$ws.Range("I60").Value = 0.9979324302179553
$ws.Range("BG60").Value = 0.96520131659323805
$ws.Range("BJ60").Value = 0.71533543277103528
$ws.Range("BO60").Value = 0.7616182333781345
$ws.Range("AV61").Value = 0.91533224960156412
$ws.Range("BK61").Value = 0.81560836267824155
$ws.Range("BI62").Value = 0.75236940967124533
$ws.Range("T63").Value = 0.69961530074854494
$ws.Range("BJ63").Value = 0.86421951268731312
$ws.Range("BO63").Value = 0.98232479015108587
$ws.Range("BJ64").Value = 0.73204840599069065
$ws.Range("BM64").Value = 0.94578101048281571
$ws.Range("BP64").Value = 0.64031575759279746
$ws.Range("BO65").Value = 0.81324198730321884
$ws.Range("AV66").Value = 0.6744888053976732
$ws.Range("BM66").Value = 0.73450245269587211
$ws.Range("BN67").Value = 0.90667183200529688
$ws.Range("BP67").Value = 0.90283765699675478
